## Add the "176. Second Highest Salary" row to the Pandas30 question table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The data lives in an Excel Table ("Table2"), currently A1:E5 with a header
# row. Growing it by one ListRow extends the table ref (and the sheet
# dimension / autofilter) from A1:E5 to A1:E6, exactly like typing into the
# row right below the table in the UI would.
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# Question / Pattern / Notes text for row 6.
$ws.Range("A6").Value = "176. Second Highest Salary"
$ws.Range("C6").Value = "Data Manipulation"
$ws.Range("D6").Value = "drop duplicates, sort descending, check if len is greater than 1, if it is get the 2nd largest value with .nlargest(2).iloc[-1], else set it as None. If none, return return pd.DataFrame({'SecondHighestSalary': [None]}). If there is one, set res_df = pd.DataFrame({'SecondHighestSalary': [second_highest]}), then return res_df."

# Difficulty = "Medium", same as row 5, so copy that cell's formatting
# (orange fill) across onto the new cell too.
$ws.Range("B6").Value = "Medium"
$ws.Range("B5").Copy() | Out-Null
$ws.Range("B6").PasteSpecial(-4122) | Out-Null

# Link cell: set the URL text, wire up the hyperlink, then pick up the
# Hyperlink cell style used by the other Link cells in the column.
$url = "https://leetcode.com/problems/second-highest-salary/solutions/3859199/pandas-an-effortless-and-simple-approach-with-comments/?envType=study-plan-v2&envId=30-days-of-pandas&lang=pythondata "
$ws.Range("E6").Value = $url
$ws.Hyperlinks.Add($ws.Range("E6"), $url) | Out-Null
$ws.Range("E5").Copy() | Out-Null
$ws.Range("E6").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Match the saved selection state.
$ws.Range("E15").Select() | Out-Null
